$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split the run "translator between your .NET objects (entities)
# and " into three runs so that the word "between" can be wrapped with
# <w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>
# markers (as Word's grammar checker would do), while leaving the visible
# text and highlight formatting unchanged.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("translator between", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null

if ($rng.Find.Found) {
    # Grow the found range out to the whole enclosing paragraph so the
    # replacement XML can be written for the complete paragraph in one shot.
    $rng.Expand(4) | Out-Null   # wdParagraph

    $paraXml = '<w:p w14:paraId="1D779692" w14:textId="2431496E" w:rsidR="00A43246" w:rsidRDefault="00A43246" w:rsidP="00A43246" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
               '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
               '<w:r w:rsidRPr="00A43246"><w:t xml:space="preserve">EF Core acts as a </w:t></w:r>' + `
               '<w:r w:rsidRPr="00A30DD7"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">translator </w:t></w:r>' + `
               '<w:proofErr w:type="gramStart"/>' + `
               '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>between</w:t></w:r>' + `
               '<w:proofErr w:type="gramEnd"/>' + `
               '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> your .NET objects (entities) and </w:t></w:r>' + `
               '<w:r w:rsidR="005404CC" w:rsidRPr="00A30DD7"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>tables</w:t></w:r>' + `
               '<w:r w:rsidRPr="00A43246"><w:t xml:space="preserve"> in a relational database. It maps classes to tables and properties to columns.</w:t></w:r>' + `
               '</w:p>'

    $pkgXml = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?>" + `
              "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" + `
              "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" + `
              "<pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>$paraXml</w:body></w:document></pkg:xmlData>" + `
              "</pkg:part></pkg:package>"

    $rng.InsertXML($pkgXml)
}

# ---------------------------------------------------------------------------
# Change 2: remove one of the three trailing empty paragraphs just before
# the section properties, leaving two.
# ---------------------------------------------------------------------------
$d2 = $word.ActiveDocument
$count = $d2.Paragraphs.Count
$trailingEmpty = @()
for ($i = $count; $i -ge 1; $i--) {
    $para = $d2.Paragraphs.Item($i)
    if ($para.Range.Text -eq "`r") {
        $trailingEmpty += $i
    } else {
        break
    }
}

if ($trailingEmpty.Count -ge 3) {
    # Delete the earliest (topmost) of the trailing empty paragraphs, as in
    # the diff, leaving the other two untouched.
    $targetIndex = ($trailingEmpty | Measure-Object -Minimum).Minimum
    $d2.Paragraphs.Item($targetIndex).Range.Delete() | Out-Null
}
